# "Generate Report for Archive"
# - Update status text "Ready for handoff" -> "In Translation" everywhere it appears
#   (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# - Narrow the "Status"/"zh-cn"/"de-de" columns (Overview cols E:F, zh-cn/de-de col C)
#   from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- 1. Replace the status text on every sheet ---------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Resize the relevant columns ---------------------------------------
# NOTE: Range.ColumnWidth (like real Excel) is quantized to whole pixels
# (character-width units are stored as pixel-count/MaxDigitWidth), so the
# closest reachable value to the target 13.4101845877511 is 13.3333333...
# (i.e. an input of ~12.5 "characters").
$targetWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").EntireColumn.ColumnWidth = $targetWidth
$overview.Range("F1").EntireColumn.ColumnWidth = $targetWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").EntireColumn.ColumnWidth = $targetWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").EntireColumn.ColumnWidth = $targetWidth
